$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target sheet is a full 4x4 Sending-cluster x Target-cluster grid (ECs/FAPs/M2/sCs),
# replacing the earlier 4x3 grid that omitted the "M2" target cluster; expression
# stats were recomputed with the revised (Dr Hou-advised) parameters.

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'ECs'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'ECs'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 5.423951666666667
$row[0,7] = 16.271855
$row[0,8] = 0.4774188439413272
$row[0,9] = 0.4774188439413271
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 5.722664999999999
$row[0,13] = 17.167995
$row[0,14] = 0.09021166427595352
$row[0,15] = 0.09021166427595351
$row[0,16] = 31.039458364525
$row[0,17] = 279.355125280725
$row[0,18] = 0.04306874846864886
$row[0,19] = 0.04306874846864884
$ws.Range("A2:T2").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'ECs'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'FAPs'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 5.423951666666667
$row[0,7] = 16.271855
$row[0,8] = 0.4774188439413272
$row[0,9] = 0.4774188439413271
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 41.286995
$row[0,13] = 123.860985
$row[0,14] = 0.6508451100847196
$row[0,15] = 0.6508451100847196
$row[0,16] = 223.9386653419083
$row[0,17] = 2015.447988077175
$row[0,18] = 0.3107257200415127
$row[0,19] = 0.3107257200415126
$ws.Range("A3:T3").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'ECs'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'M2'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 5.423951666666667
$row[0,7] = 16.271855
$row[0,8] = 0.4774188439413272
$row[0,9] = 0.4774188439413271
$row[0,10] = 1
$row[0,11] = 0.3333333333333333
$row[0,12] = 0.06212466666666667
$row[0,13] = 0.186374
$row[0,14] = 0.0009793286122093212
$row[0,15] = 0.000979328612209321
$row[0,16] = 0.3369611893077778
$row[0,17] = 3.03265070377
$row[0,18] = 0.0004675499338796385
$row[0,19] = 0.0004675499338796383
$ws.Range("A4:T4").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'ECs'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'sCs'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 5.423951666666667
$row[0,7] = 16.271855
$row[0,8] = 0.4774188439413272
$row[0,9] = 0.4774188439413271
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 16.36419166666667
$row[0,13] = 49.092575
$row[0,14] = 0.2579638970271176
$row[0,15] = 0.2579638970271176
$row[0,16] = 88.75858466406945
$row[0,17] = 798.827261976625
$row[0,18] = 0.1231568254972861
$row[0,19] = 0.123156825497286
$ws.Range("A5:T5").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'FAPs'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'ECs'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 1.583504333333333
$row[0,7] = 4.750513
$row[0,8] = 0.1393808158066948
$row[0,9] = 0.1393808158066948
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 5.722664999999999
$row[0,13] = 17.167995
$row[0,14] = 0.09021166427595352
$row[0,15] = 0.09021166427595351
$row[0,16] = 9.061864825714999
$row[0,17] = 81.55678343143498
$row[0,18] = 0.01257377536206207
$row[0,19] = 0.01257377536206206
$ws.Range("A6:T6").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'FAPs'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'FAPs'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 1.583504333333333
$row[0,7] = 4.750513
$row[0,8] = 0.1393808158066948
$row[0,9] = 0.1393808158066948
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 41.286995
$row[0,13] = 123.860985
$row[0,14] = 0.6508451100847196
$row[0,15] = 0.6508451100847196
$row[0,16] = 65.37813549281165
$row[0,17] = 588.403219435305
$row[0,18] = 0.0907153224074063
$row[0,19] = 0.09071532240740629
$ws.Range("A7:T7").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'FAPs'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'M2'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 1.583504333333333
$row[0,7] = 4.750513
$row[0,8] = 0.1393808158066948
$row[0,9] = 0.1393808158066948
$row[0,10] = 1
$row[0,11] = 0.3333333333333333
$row[0,12] = 0.06212466666666667
$row[0,13] = 0.186374
$row[0,14] = 0.0009793286122093212
$row[0,15] = 0.000979328612209321
$row[0,16] = 0.09837467887355555
$row[0,17] = 0.885372109862
$row[0,18] = 0.0001364996209125735
$row[0,19] = 0.0001364996209125734
$ws.Range("A8:T8").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'FAPs'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'sCs'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 1.583504333333333
$row[0,7] = 4.750513
$row[0,8] = 0.1393808158066948
$row[0,9] = 0.1393808158066948
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 16.36419166666667
$row[0,13] = 49.092575
$row[0,14] = 0.2579638970271176
$row[0,15] = 0.2579638970271176
$row[0,16] = 25.91276841566389
$row[0,17] = 233.214915740975
$row[0,18] = 0.03595521841631386
$row[0,19] = 0.03595521841631386
$ws.Range("A9:T9").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'M2'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'ECs'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 0.6660723333333333
$row[0,7] = 1.998217
$row[0,8] = 0.05862800830537802
$row[0,9] = 0.05862800830537802
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 5.722664999999999
$row[0,13] = 17.167995
$row[0,14] = 0.09021166427595352
$row[0,15] = 0.09021166427595351
$row[0,16] = 3.811708829434999
$row[0,17] = 34.30537946491499
$row[0,18] = 0.005288930202412577
$row[0,19] = 0.005288930202412576
$ws.Range("A10:T10").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'M2'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'FAPs'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 0.6660723333333333
$row[0,7] = 1.998217
$row[0,8] = 0.05862800830537802
$row[0,9] = 0.05862800830537802
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 41.286995
$row[0,13] = 123.860985
$row[0,14] = 0.6508451100847196
$row[0,15] = 0.6508451100847196
$row[0,16] = 27.50012509597166
$row[0,17] = 247.501125863745
$row[0,18] = 0.03815775251956161
$row[0,19] = 0.03815775251956161
$ws.Range("A11:T11").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'M2'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'M2'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 0.6660723333333333
$row[0,7] = 1.998217
$row[0,8] = 0.05862800830537802
$row[0,9] = 0.05862800830537802
$row[0,10] = 1
$row[0,11] = 0.3333333333333333
$row[0,12] = 0.06212466666666667
$row[0,13] = 0.186374
$row[0,14] = 0.0009793286122093212
$row[0,15] = 0.000979328612209321
$row[0,16] = 0.04137952168422222
$row[0,17] = 0.372415695158
$row[0,18] = 0.000057416086010302416764439215
$row[0,19] = 0.000057416086010302403211912059
$ws.Range("A12:T12").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'M2'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'sCs'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 0.6660723333333333
$row[0,7] = 1.998217
$row[0,8] = 0.05862800830537802
$row[0,9] = 0.05862800830537802
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 16.36419166666667
$row[0,13] = 49.092575
$row[0,14] = 0.2579638970271176
$row[0,15] = 0.2579638970271176
$row[0,16] = 10.89973532653056
$row[0,17] = 98.097617938775
$row[0,18] = 0.01512390949739353
$row[0,19] = 0.01512390949739353
$ws.Range("A13:T13").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'sCs'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'ECs'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 3.687463666666666
$row[0,7] = 11.062391
$row[0,8] = 0.3245723319466
$row[0,9] = 0.3245723319466
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 5.722664999999999
$row[0,13] = 17.167995
$row[0,14] = 0.09021166427595352
$row[0,15] = 0.09021166427595351
$row[0,16] = 21.102119264005
$row[0,17] = 189.919073376045
$row[0,18] = 0.02928021024283002
$row[0,19] = 0.02928021024283001
$ws.Range("A14:T14").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'sCs'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'FAPs'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 3.687463666666666
$row[0,7] = 11.062391
$row[0,8] = 0.3245723319466
$row[0,9] = 0.3245723319466
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 41.286995
$row[0,13] = 123.860985
$row[0,14] = 0.6508451100847196
$row[0,15] = 0.6508451100847196
$row[0,16] = 152.2442939683483
$row[0,17] = 1370.198645715135
$row[0,18] = 0.211246315116239
$row[0,19] = 0.211246315116239
$ws.Range("A15:T15").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'sCs'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'M2'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 3.687463666666666
$row[0,7] = 11.062391
$row[0,8] = 0.3245723319466
$row[0,9] = 0.3245723319466
$row[0,10] = 1
$row[0,11] = 0.3333333333333333
$row[0,12] = 0.06212466666666667
$row[0,13] = 0.186374
$row[0,14] = 0.0009793286122093212
$row[0,15] = 0.000979328612209321
$row[0,16] = 0.2290824511371111
$row[0,17] = 2.061742060234
$row[0,18] = 0.0003178629714068069
$row[0,19] = 0.0003178629714068068
$ws.Range("A16:T16").Value = $row

$row = New-Object 'object[,]' 1,20
$row[0,0] = 'sCs'
$row[0,1] = 'Bmp2'
$row[0,2] = 'Bmpr1a'
$row[0,3] = 'sCs'
$row[0,4] = 3
$row[0,5] = 1
$row[0,6] = 3.687463666666666
$row[0,7] = 11.062391
$row[0,8] = 0.3245723319466
$row[0,9] = 0.3245723319466
$row[0,10] = 3
$row[0,11] = 1
$row[0,12] = 16.36419166666667
$row[0,13] = 49.092575
$row[0,14] = 0.2579638970271176
$row[0,15] = 0.2579638970271176
$row[0,16] = 60.34236220520278
$row[0,17] = 543.081259846825
$row[0,18] = 0.08372794361612414
$row[0,19] = 0.08372794361612414
$ws.Range("A17:T17").Value = $row

